# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (with fund holding data) between the
# "总计" (totals) sheet and the existing "2022-Q2" sheet, and updates the
# "总计" sheet so that:
#   - row 2 now reports the new 2022-Q3 totals (6 holdings, 0.12 亿元)
#   - a new row 3 preserves the old 2022-Q2 totals (8 holdings, 1.2 亿元)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 总计 (totals sheet; the existing "2022-Q2"
                                 # fund-holding sheet stays at index 2 until
                                 # the new sheet is inserted after $ws1 below)

# ---------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet: shift the old 2022-Q2 row down to
#    row 3 and put the new 2022-Q3 totals in row 2.
# ---------------------------------------------------------------------
$oldB2 = $ws1.Range("B2").Value2
$oldC2 = $ws1.Range("C2").Value2
$oldD2 = $ws1.Range("D2").Value2

# Duplicate row 2's formatting into row 3 (keeps the A-column header style).
$ws1.Range("A2:D2").Copy($ws1.Range("A3:D3"))

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = $oldB2
$ws1.Range("C3").Value = $oldC2
$ws1.Range("D3").Value = $oldD2

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 6
$ws1.Range("D2").Value = 0.12

# ---------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" (so it lands
#    before the existing "2022-Q2" sheet) and fill it with data.
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "2022-Q3"

# Bring over header / index-column styling from "总计" (s=2 style class).
$ws1.Range("B1:D1").Copy($newSheet.Range("B1:D1"))
$ws1.Range("B1:D1").Copy($newSheet.Range("E1:G1"))
$ws1.Range("B1").Copy($newSheet.Range("H1"))
$ws1.Range("A2").Copy($newSheet.Range("A2:A7"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G hold text (numeric-looking strings must keep their leading
# zeros / trailing zeros, e.g. fund code "004895" or ratio "0.0920"), so
# they're written with a leading apostrophe to force text, just like a
# user typing '004895 into a cell. H (and the lone G7) are real numbers.
$data = @(
    @(0, "'004895", "华商鑫安灵活配置混合",       "'2.11", "'92.54", "'4.36", "'0.0920", 3),
    @(1, "'004258", "国寿安保稳嘉混合A",           "'2.32", "'20.12", "'0.73", "'0.0169", 9),
    @(2, "'007533", "格林创新成长混合A",           "'0.05", "'88.63", "'7.64", "'0.0038", 2),
    @(3, "'007534", "格林创新成长混合C",           "'0.04", "'88.63", "'7.64", "'0.0031", 2),
    @(4, "'970083", "东海证券海盈6个月持有期混合", "'0.12", "'36.37", "'2.49", "'0.0030", 4),
    @(5, "'004259", "国寿安保稳嘉混合C",           "'0.00", "'20.12", "'0.73", 0,         9)
)

$row = 2
foreach ($rec in $data) {
    $newSheet.Range("A" + $row).Value = $rec[0]
    $newSheet.Range("B" + $row).Value = $rec[1]
    $newSheet.Range("C" + $row).Value = $rec[2]
    $newSheet.Range("D" + $row).Value = $rec[3]
    $newSheet.Range("E" + $row).Value = $rec[4]
    $newSheet.Range("F" + $row).Value = $rec[5]
    $newSheet.Range("G" + $row).Value = $rec[6]
    $newSheet.Range("H" + $row).Value = $rec[7]
    # Drop the quote-prefix formatting the apostrophe-forced-text entries
    # above picked up, so only the value's type (text) changed, not its look.
    $newSheet.Range("B" + $row + ":G" + $row).Style = "Normal"
    $row = $row + 1
}
